$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 ("R") updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 227
$wsOff.Range("C3").Value = 143
$wsOff.Range("D3").Value = 45
$wsOff.Range("E3").Value = 17

# DEF sheet - row 3 ("R") updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 232
$wsDef.Range("C3").Value = 174
$wsDef.Range("D3").Value = 55
$wsDef.Range("E3").Value = 30
